# Updated table: compare glmnet vs svm - added test errors
# (mirrors the edits made to predictivemodelingVY.rmd)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the " ROC " column headers (row 2) to "AUC" for both the glmnet
# and SVM blocks.
$ws.Range("B2").Value = "AUC"
$ws.Range("E2").Value = "AUC"

# The Testing row used to hold raw ROC/Sens/Spec numbers for SVM only
# (E4:G4). Replace that with the new AUC +/- error summary values for
# both models, reported as text in B4 (glmnet) and E4 (SVM).
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null
$ws.Range("G4").Value = $null

$ws.Range("B4").Value = "0.977+-0.024"
$ws.Range("E4").Value = "0.978+-0.05"

# Drop the leftover helper table (K9:N9 labels, J10:J12 values) that is
# no longer needed.
$ws.Range("K9:N9").Clear()
$ws.Range("J10:J12").Clear()

# Match the final saved selection state.
$ws.Range("C9").Select()
